$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    3  = -1
    4  = 1
    5  = 2
    6  = -3
    8  = -4
    9  = 2
    10 = -2
    11 = -3
    12 = -3
    13 = -1
    15 = 3
    16 = -4
    18 = -4
    19 = 2
    20 = -3
    21 = 0
    22 = -5
    23 = 6
    24 = -3
    25 = 2
    26 = 1
    27 = -5
    28 = 2
    29 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
